# Updated symbol list on Sat Dec 24 05:59:02 UTC 2022 with GitHub Actions
#
# Refreshes the "Price" (column D) and a couple of the concatenated
# "Worstin24h"-tagged helper strings (column E) on Sheet1, matching a
# fresh pull of the coin-ranking data. Price values are numeric-looking
# text (the sheet stores them as strings), so each is written with a
# leading apostrophe to force text entry and then ClearFormats() is used
# to drop the quote-prefix styling Excel applies, keeping cell style
# untouched (same as before the edit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'245.43"
$ws.Range("D2").ClearFormats()
$ws.Range("D4").Value = "'5.323"
$ws.Range("D4").ClearFormats()
$ws.Range("D5").Value = "'0.05963"
$ws.Range("D5").ClearFormats()
$ws.Range("D6").Value = "'3.396"
$ws.Range("D6").ClearFormats()
$ws.Range("D7").Value = "'6.390"
$ws.Range("D7").ClearFormats()
$ws.Range("D8").Value = "'0.8100"
$ws.Range("D8").ClearFormats()
$ws.Range("D9").Value = "'0.9634"
$ws.Range("D9").ClearFormats()
$ws.Range("D10").Value = "'0.1429"
$ws.Range("D10").ClearFormats()
$ws.Range("D11").Value = "'0.07403"
$ws.Range("D11").ClearFormats()
$ws.Range("D12").Value = "'0.03422"
$ws.Range("D12").ClearFormats()
$ws.Range("D14").Value = "'0.09403"
$ws.Range("D14").ClearFormats()
$ws.Range("D15").Value = "'3.997"
$ws.Range("D15").ClearFormats()
$ws.Range("D16").Value = "'0.001596"
$ws.Range("D16").ClearFormats()
$ws.Range("D17").Value = "'0.04813"
$ws.Range("D17").ClearFormats()
$ws.Range("D18").Value = "'0.0005913"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "17OneONEWorstin24h"
$ws.Range("D19").Value = "'0.006230"
$ws.Range("D19").ClearFormats()
$ws.Range("D20").Value = "'0.004141"
$ws.Range("D20").ClearFormats()
$ws.Range("D21").Value = "'0.0009879"
$ws.Range("D21").ClearFormats()
$ws.Range("D22").Value = "'0.00009707"
$ws.Range("D22").ClearFormats()
$ws.Range("D23").Value = "'3.739"
$ws.Range("D23").ClearFormats()
$ws.Range("D26").Value = "'0.1332"
$ws.Range("D26").ClearFormats()
$ws.Range("D27").Value = "'0.0002462"
$ws.Range("D27").ClearFormats()
$ws.Range("D40").Value = "'0.03911"
$ws.Range("D40").ClearFormats()
$ws.Range("D41").Value = "'0.006455"
$ws.Range("D41").ClearFormats()
$ws.Range("D42").Value = "'0.1072"
$ws.Range("D42").ClearFormats()
$ws.Range("D43").Value = "'0.003002"
$ws.Range("D43").ClearFormats()
$ws.Range("D44").Value = "'0.005349"
$ws.Range("D44").ClearFormats()
$ws.Range("D45").Value = "'0.00005323"
$ws.Range("D45").ClearFormats()
$ws.Range("D47").Value = "'0.8503"
$ws.Range("D47").ClearFormats()
$ws.Range("D48").Value = "'0.03830"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "47BOLOBOLO"
$ws.Range("D49").Value = "'0.00002101"
$ws.Range("D49").ClearFormats()
